# Reorders the three "Aydin Biber" user-story bullets:
#   voorwerpen -> springen  (gains the _GoBack bookmark)
#   vijanden   -> voorwerpen
#   springen   -> vijanden
# i.e. the "springen" bullet (with its _GoBack bookmark) moves to the top,
# and the other two shift down one slot.

$d = $word.ActiveDocument

# The _GoBack bookmark currently sits at the end of the "springen" bullet.
# Drop it now; it gets re-created at the correct (new) location below.
$bm = $d.Bookmarks("_GoBack")
$bmDeleted = $bm.Delete()

# Rotate the three bullet texts via a temporary placeholder so none of the
# three Find/Replace calls can accidentally match text just written by a
# previous one.
$f1 = $d.Content.Find.Execute( `
    "Als speler wil ik voorwerpen kunnen oppakken om mijn karakter sterker te maken", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "PLACEHOLDER_ROTATE_TEXT", 2)

$f2 = $d.Content.Find.Execute( `
    "Als een speler wil ik vijanden kunnen zien aanvallen om deze te ontwijken", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Als speler wil ik voorwerpen kunnen oppakken om mijn karakter sterker te maken", 2)

$f3 = $d.Content.Find.Execute( `
    "Als een speler wil ik springen zodat ik op platformen kan komen", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Als een speler wil ik vijanden kunnen zien aanvallen om deze te ontwijken", 2)

$f4 = $d.Content.Find.Execute( `
    "PLACEHOLDER_ROTATE_TEXT", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Als een speler wil ik springen zodat ik op platformen kan komen", 2)

# Re-add the _GoBack bookmark at the end of the "springen" bullet's text,
# which is now the first bullet in the list.
#
# A Range built directly at "end of this paragraph's text" (i.e. a
# zero-length Range right before the paragraph mark) is mis-resolved by
# Bookmarks.Add in this runtime. Work around it by inserting a one-
# character placeholder there, bookmarking that one character, and then
# deleting the placeholder - the bookmark collapses onto the correct point
# once the bracketed placeholder text disappears.
$paras = $d.Paragraphs
$targetPara = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text.Contains("springen zodat ik op platformen kan komen")) {
        $targetPara = $cand
    }
}

$targetPara.Range.InsertAfter("X")
$xPos = $targetPara.Range.End - 2
$placeholderRange = $d.Range($xPos, $xPos + 1)
$bmAdded = $d.Bookmarks.Add("_GoBack", $placeholderRange)
$deleteRange = $d.Range($xPos, $xPos + 1)
$deleteRange.Delete()
